$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove duplicate cell values in column E (row 4 and row 8)
$ws.Range("E4").ClearContents()
$ws.Range("E8").ClearContents()

# Update the selected cell/range shown in the sheet view
$ws.Range("G10").Select()
